# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    (Overview sheet mirrors this through its zh-cn / de-de summary columns)
#  - Each language sheet's "Latest Handback DateTime" is refreshed
#  - The stale "handback file is not the latest" Error Detail is cleared
#    now that the handback is in sync

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-01 15:11:45"
$wsZhCn.Range("P2").Value = ""

# ---- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-01 15:11:54"
$wsDeDe.Range("P2").Value = ""

# ---- Column widths: Status columns grew, now-empty Error Detail columns shrank ----
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

$wsZhCn.Columns.Item(3).ColumnWidth = 29.14
$wsZhCn.Columns.Item(16).ColumnWidth = 12.91

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14
$wsDeDe.Columns.Item(16).ColumnWidth = 12.91
